$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell $ws.Cells.Item(2, 4) "68.356.16"
$ws.Cells.Item(2, 5).Value = "  +1.53%  "
Set-TextCell $ws.Cells.Item(3, 4) "2.640.30"
$ws.Cells.Item(3, 5).Value = "  +1.16%  "
$ws.Cells.Item(4, 5).Value = "  +0.03%  "
Set-TextCell $ws.Cells.Item(5, 4) "599.24"
$ws.Cells.Item(5, 5).Value = "  +1.30%  "
Set-TextCell $ws.Cells.Item(6, 4) "154.15"
$ws.Cells.Item(6, 5).Value = "  +2.47%  "
$ws.Cells.Item(7, 5).Value = "  -0.01%  "
$ws.Cells.Item(8, 5).Value = "  +0.00%  "
Set-TextCell $ws.Cells.Item(9, 4) "2.639.11"
$ws.Cells.Item(9, 5).Value = "  +1.19%  "
$ws.Cells.Item(10, 5).Value = "  +8.17%  "
$ws.Cells.Item(11, 5).Value = "  -0.56%  "
$ws.Cells.Item(12, 5).Value = "  +1.08%  "
Set-TextCell $ws.Cells.Item(13, 4) "0.349"
$ws.Cells.Item(13, 5).Value = "  +1.60%  "
Set-TextCell $ws.Cells.Item(14, 4) "27.93"
$ws.Cells.Item(14, 5).Value = "  +2.53%  "
$ws.Cells.Item(15, 5).Value = "  +3.20%  "
$ws.Cells.Item(16, 5).Value = "  +1.54%  "
Set-TextCell $ws.Cells.Item(17, 4) "68.283.99"
$ws.Cells.Item(17, 5).Value = "  +1.28%  "
Set-TextCell $ws.Cells.Item(18, 4) "2.646.25"
$ws.Cells.Item(18, 5).Value = "  +1.37%  "
$ws.Cells.Item(19, 5).Value = "  +3.58%  "
Set-TextCell $ws.Cells.Item(20, 4) "366.18"
$ws.Cells.Item(20, 5).Value = "  -1.86%  "
Set-TextCell $ws.Cells.Item(21, 4) "7.40"
$ws.Cells.Item(21, 5).Value = "  +0.40%  "
$ws.Cells.Item(22, 5).Value = "  -0.59%  "
$ws.Cells.Item(23, 5).Value = "  +0.15%  "
$ws.Cells.Item(24, 5).Value = "  +3.11%  "
Set-TextCell $ws.Cells.Item(25, 4) "73.68"
$ws.Cells.Item(25, 5).Value = "  -0.08%  "
Set-TextCell $ws.Cells.Item(26, 4) "1.00"
$ws.Cells.Item(26, 5).Value = "  +0.00%  "
Set-TextCell $ws.Cells.Item(27, 4) "9.99"
$ws.Cells.Item(27, 5).Value = "  +0.64%  "
Set-TextCell $ws.Cells.Item(28, 4) "2.773.37"
$ws.Cells.Item(29, 5).Value = "  +5.76%  "
Set-TextCell $ws.Cells.Item(31, 4) "573.61"
$ws.Cells.Item(31, 5).Value = "  -0.70%  "
$ws.Cells.Item(32, 5).Value = "  +4.30%  "
Set-TextCell $ws.Cells.Item(33, 4) "7.99"
$ws.Cells.Item(33, 5).Value = "  +4.37%  "
$ws.Cells.Item(34, 5).Value = "  +2.52%  "
Set-TextCell $ws.Cells.Item(35, 4) "0.130"
$ws.Cells.Item(35, 5).Value = "  +2.72%  "
Set-TextCell $ws.Cells.Item(36, 4) "1.00"
$ws.Cells.Item(36, 5).Value = "  +0.03%  "
$ws.Cells.Item(37, 5).Value = "  +3.33%  "
Set-TextCell $ws.Cells.Item(38, 4) "160.25"
$ws.Cells.Item(38, 5).Value = "  +1.65%  "
Set-TextCell $ws.Cells.Item(39, 4) "19.25"
$ws.Cells.Item(39, 5).Value = "  +1.05%  "
$ws.Cells.Item(40, 5).Value = "  +3.85%  "
$ws.Cells.Item(41, 5).Value = "  +0.88%  "
$ws.Cells.Item(42, 5).Value = "  +2.52%  "
$ws.Cells.Item(43, 5).Value = "  +3.56%  "
$ws.Cells.Item(44, 5).Value = "  +2.88%  "
Set-TextCell $ws.Cells.Item(45, 4) "0.0₆0319"
$ws.Cells.Item(45, 5).Value = "  +12.50%  "
$ws.Cells.Item(46, 5).Value = "  +0.05%  "
Set-TextCell $ws.Cells.Item(47, 4) "40.51"
$ws.Cells.Item(47, 5).Value = "  -0.35%  "
Set-TextCell $ws.Cells.Item(48, 4) "157.00"
$ws.Cells.Item(48, 5).Value = "  +2.60%  "
Set-TextCell $ws.Cells.Item(49, 4) "3.74"
$ws.Cells.Item(49, 5).Value = "  +0.86%  "
Set-TextCell $ws.Cells.Item(50, 4) "1.71"
$ws.Cells.Item(50, 5).Value = "  +2.02%  "
Set-TextCell $ws.Cells.Item(51, 4) "21.86"
$ws.Cells.Item(51, 5).Value = "  +2.39%  "
